$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B:F, rows 2-25
$dataBF = New-Object "object[,]" 24,5
$dataBF[0,0] = 1.02
$dataBF[0,1] = 1.035313136245032
$dataBF[0,2] = 1.045885087409457
$dataBF[0,3] = 1.052805140006511
$dataBF[0,4] = 1.0583272519715
$dataBF[1,0] = 1.02
$dataBF[1,1] = 1.036097815517028
$dataBF[1,2] = 1.046532055325228
$dataBF[1,3] = 1.053680937752309
$dataBF[1,4] = 1.059172040103021
$dataBF[2,0] = 1.02
$dataBF[2,1] = 1.036606015148782
$dataBF[2,2] = 1.04695106004817
$dataBF[2,3] = 1.054248984895267
$dataBF[2,4] = 1.059719698591816
$dataBF[3,0] = 1.02
$dataBF[3,1] = 1.036819770597775
$dataBF[3,2] = 1.047127297153307
$dataBF[3,3] = 1.054488112567262
$dataBF[3,4] = 1.059950177415448
$dataBF[4,0] = 1.02
$dataBF[4,1] = 1.036855667380329
$dataBF[4,2] = 1.047156893217677
$dataBF[4,3] = 1.054528281902857
$dataBF[4,4] = 1.059988890027838
$dataBF[5,0] = 1.02
$dataBF[5,1] = 1.036608870936071
$dataBF[5,2] = 1.046953414595403
$dataBF[5,3] = 1.054252178873394
$dataBF[5,4] = 1.059722777308023
$dataBF[6,0] = 1.02
$dataBF[6,1] = 1.035578226064459
$dataBF[6,2] = 1.046103654868505
$dataBF[6,3] = 1.053100839800707
$dataBF[6,4] = 1.058612538848348
$dataBF[7,0] = 1.02
$dataBF[7,1] = 1.033765694863568
$dataBF[7,2] = 1.044609208921284
$dataBF[7,3] = 1.051082438263127
$dataBF[7,4] = 1.056664087244534
$dataBF[8,0] = 1.02
$dataBF[8,1] = 1.032559861566806
$dataBF[8,2] = 1.043614996903641
$dataBF[8,3] = 1.049743948032356
$dataBF[8,4] = 1.055370569120651
$dataBF[9,0] = 1.02
$dataBF[9,1] = 1.032038341935415
$dataBF[9,2] = 1.043185009769091
$dataBF[9,3] = 1.049166077888386
$dataBF[9,4] = 1.054811778918844
$dataBF[10,0] = 1.02
$dataBF[10,1] = 1.031844720325589
$dataBF[10,2] = 1.043025372373506
$dataBF[10,3] = 1.048951689313235
$dataBF[10,4] = 1.054604418436957
$dataBF[11,0] = 1.02
$dataBF[11,1] = 1.031886248536524
$dataBF[11,2] = 1.043059611533235
$dataBF[11,3] = 1.048997664677902
$dataBF[11,4] = 1.05464888892314
$dataBF[12,0] = 1.02
$dataBF[12,1] = 1.032022335185103
$dataBF[12,2] = 1.043171812470409
$dataBF[12,3] = 1.049148351184769
$dataBF[12,4] = 1.054794634363923
$dataBF[13,0] = 1.02
$dataBF[13,1] = 1.032106195168774
$dataBF[13,2] = 1.04324095369109
$dataBF[13,3] = 1.049241228378218
$dataBF[13,4] = 1.054884459380354
$dataBF[14,0] = 1.02
$dataBF[14,1] = 1.03259448619442
$dataBF[14,2] = 1.043643544706002
$dataBF[14,3] = 1.049782335460757
$dataBF[14,4] = 1.055407681993549
$dataBF[15,0] = 1.02
$dataBF[15,1] = 1.032900943901498
$dataBF[15,2] = 1.043896218136554
$dataBF[15,3] = 1.050122215497885
$dataBF[15,4] = 1.055736238311642
$dataBF[16,0] = 1.02
$dataBF[16,1] = 1.033079754636773
$dataBF[16,2] = 1.044043647668064
$dataBF[16,3] = 1.050320626144528
$dataBF[16,4] = 1.055928006076487
$dataBF[17,0] = 1.02
$dataBF[17,1] = 1.033140734454989
$dataBF[17,2] = 1.04409392564977
$dataBF[17,3] = 1.050388306867781
$dataBF[17,4] = 1.055993415310878
$dataBF[18,0] = 1.02
$dataBF[18,1] = 1.032868057760218
$dataBF[18,2] = 1.043869103546189
$dataBF[18,3] = 1.050085732570551
$dataBF[18,4] = 1.055700974237061
$dataBF[19,0] = 1.02
$dataBF[19,1] = 1.031982258450123
$dataBF[19,2] = 1.043138769895992
$dataBF[19,3] = 1.049103970622854
$dataBF[19,4] = 1.054751710452491
$dataBF[20,0] = 1.02
$dataBF[20,1] = 1.031425865992095
$dataBF[20,2] = 1.042680037968061
$dataBF[20,3] = 1.048488192335654
$dataBF[20,4] = 1.054156023442418
$dataBF[21,0] = 1.02
$dataBF[21,1] = 1.031720767948378
$dataBF[21,2] = 1.042923176424239
$dataBF[21,3] = 1.048814485743659
$dataBF[21,4] = 1.054471698452039
$dataBF[22,0] = 1.02
$dataBF[22,1] = 1.032882917409376
$dataBF[22,2] = 1.04388135531006
$dataBF[22,3] = 1.050102217127397
$dataBF[22,4] = 1.05571690816407
$dataBF[23,0] = 1.02
$dataBF[23,1] = 1.034233840826107
$dataBF[23,2] = 1.044995199452532
$dataBF[23,3] = 1.051602998540882
$dataBF[23,4] = 1.057166856769974

# New values for columns I:N, rows 2-25
$dataIN = New-Object "object[,]" 24,6
$dataIN[0,0] = 1.041396201793307
$dataIN[0,1] = 1.040427594019254
$dataIN[0,2] = 1.048652073421774
$dataIN[0,3] = 1.055552857427916
$dataIN[0,4] = 1.061059789117118
$dataIN[0,5] = 1.017443140021282
$dataIN[1,0] = 1.041595204610035
$dataIN[1,1] = 1.040856601863624
$dataIN[1,2] = 1.049111000415069
$dataIN[1,3] = 1.056241412644247
$dataIN[1,4] = 1.061718510450427
$dataIN[1,5] = 1.017585527122665
$dataIN[2,0] = 1.041722832677424
$dataIN[2,1] = 1.041133971317666
$dataIN[2,2] = 1.049407642483705
$dataIN[2,3] = 1.056687659201148
$dataIN[2,4] = 1.06214512234548
$dataIN[2,5] = 1.017677570288323
$dataIN[3,0] = 1.041776213751466
$dataIN[3,1] = 1.041250521789274
$dataIN[3,2] = 1.049532274176619
$dataIN[3,3] = 1.056875428476441
$dataIN[3,4] = 1.062324558135666
$dataIN[3,5] = 1.017716243058534
$dataIN[4,0] = 1.041785160590765
$dataIN[4,1] = 1.041270087827408
$dataIN[4,2] = 1.049553195852678
$dataIN[4,3] = 1.056906965522775
$dataIN[4,4] = 1.062354691315233
$dataIN[4,5] = 1.017722735073672
$dataIN[5,0] = 1.041723547034698
$dataIN[5,1] = 1.041135528890967
$dataIN[5,2] = 1.049409308121317
$dataIN[5,3] = 1.056690167528363
$dataIN[5,4] = 1.062147519631008
$dataIN[5,5] = 1.017678087123458
$dataIN[6,0] = 1.041463691283664
$dataIN[6,1] = 1.040572625555898
$dataIN[6,2] = 1.048807234220372
$dataIN[6,3] = 1.055785410875086
$dataIN[6,4] = 1.061282328400967
$dataIN[6,5] = 1.017491278949368
$dataIN[7,0] = 1.040997097553789
$dataIN[7,1] = 1.039579031357058
$dataIN[7,2] = 1.047743955768867
$dataIN[7,3] = 1.05419658762216
$dataIN[7,4] = 1.059760697530194
$dataIN[7,5] = 1.017161422864993
$dataIN[8,0] = 1.040680238838056
$dataIN[8,1] = 1.038915574435501
$dataIN[8,2] = 1.047033606815372
$dataIN[8,3] = 1.053141149247747
$dataIN[8,4] = 1.058748357175249
$dataIN[8,5] = 1.016941090154963
$dataIN[9,0] = 1.040541671156482
$dataIN[9,1] = 1.038628054081653
$dataIN[9,2] = 1.046725680083826
$dataIN[9,3] = 1.052685049155085
$dataIN[9,4] = 1.058310516942593
$dataIN[9,5] = 1.016845587468678
$dataIN[10,0] = 1.040489996580668
$dataIN[10,1] = 1.038521221530853
$dataIN[10,2] = 1.046611252541205
$dataIN[10,3] = 1.052515771902127
$dataIN[10,4] = 1.05814796186822
$dataIN[10,5] = 1.016810099360258
$dataIN[11,0] = 1.040501090187522
$dataIN[11,1] = 1.038544139025573
$dataIN[11,2] = 1.046635799870941
$dataIN[11,3] = 1.052552076147612
$dataIN[11,4] = 1.058182826913914
$dataIN[11,5] = 1.016817712313622
$dataIN[12,0] = 1.040537403884236
$dataIN[12,1] = 1.038619223970182
$dataIN[12,2] = 1.046716222476287
$dataIN[12,3] = 1.052671053806961
$dataIN[12,4] = 1.058297078482728
$dataIN[12,5] = 1.016842654296851
$dataIN[13,0] = 1.040559750894993
$dataIN[13,1] = 1.038665481722959
$dataIN[13,2] = 1.046765766934144
$dataIN[13,3] = 1.052744378288985
$dataIN[13,4] = 1.05836748306878
$dataIN[13,5] = 1.016858020012926
$dataIN[14,0] = 1.040689406413713
$dataIN[14,1] = 1.038934651281882
$dataIN[14,2] = 1.047054035832805
$dataIN[14,3] = 1.053171438458548
$dataIN[14,4] = 1.058777426074935
$dataIN[14,5] = 1.016947426344058
$dataIN[15,0] = 1.040770370820453
$dataIN[15,1] = 1.039103431086586
$dataIN[15,2] = 1.047234768989697
$dataIN[15,3] = 1.053439567275638
$dataIN[15,4] = 1.059034710177438
$dataIN[15,5] = 1.017003482869423
$dataIN[16,0] = 1.040817464157238
$dataIN[16,1] = 1.039201854335778
$dataIN[16,2] = 1.047340154593433
$dataIN[16,3] = 1.053596050131726
$dataIN[16,4] = 1.059184828718708
$dataIN[16,5] = 1.017036170287396
$dataIN[17,0] = 1.040833499382123
$dataIN[17,1] = 1.03923541015508
$dataIN[17,2] = 1.047376082723432
$dataIN[17,3] = 1.053649421592794
$dataIN[17,4] = 1.059236023534859
$dataIN[17,5] = 1.017047314243583
$dataIN[18,0] = 1.040761697735908
$dataIN[18,1] = 1.039085324997
$dataIN[18,2] = 1.04721538143465
$dataIN[18,3] = 1.053410790492051
$dataIN[18,4] = 1.059007100939691
$dataIN[18,5] = 1.016997469504379
$dataIN[19,0] = 1.040526716040632
$dataIN[19,1] = 1.038597114273109
$dataIN[19,2] = 1.04669254138836
$dataIN[19,3] = 1.052636014011555
$dataIN[19,4] = 1.058263432068812
$dataIN[19,5] = 1.016835309891764
$dataIN[20,0] = 1.040377791713175
$dataIN[20,1] = 1.038289956722888
$dataIN[20,2] = 1.046363523399599
$dataIN[20,3] = 1.052149684175093
$dataIN[20,4] = 1.057796311689608
$dataIN[20,5] = 1.016733272008127
$dataIN[21,0] = 1.040456851082996
$dataIN[21,1] = 1.038452805305327
$dataIN[21,2] = 1.046537968870904
$dataIN[21,3] = 1.052407420164735
$dataIN[21,4] = 1.058043897466387
$dataIN[21,5] = 1.016787371818179
$dataIN[22,0] = 1.040765617136922
$dataIN[22,1] = 1.039093506431216
$dataIN[22,2] = 1.04722414193835
$dataIN[22,3] = 1.053423793210629
$dataIN[22,4] = 1.059019576214375
$dataIN[22,5] = 1.017000186714315
$dataIN[23,0] = 1.041118748476948
$dataIN[23,1] = 1.039836091601605
$dataIN[23,2] = 1.048019108498691
$dataIN[23,3] = 1.054606677876569
$dataIN[23,4] = 1.060153716184206
$dataIN[23,5] = 1.017246776024066

$ws.Range("B2:F25").Value = $dataBF
$ws.Range("I2:N25").Value = $dataIN
